$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.495.40'
$ws.Range('E2').Value = '  -0.89%  '

# Row 3
$ws.Range('D3').Value = '1.896.18'
$ws.Range('E3').Value = '  -0.79%  '

# Row 4
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.61'
$ws.Range('E5').Value = '  -0.11%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9987'

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4849'
$ws.Range('E7').Value = '  -1.46%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2900'
$ws.Range('E8').Value = '  -1.99%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06636'
$ws.Range('E9').Value = '  -1.54%  '

# Row 10
$ws.Range('D10').Value = '1.904.93'
$ws.Range('E10').Value = '  +0.02%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.02'
$ws.Range('E11').Value = '  -0.22%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07399'
$ws.Range('E12').Value = '  +0.53%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.193'
$ws.Range('E13').Value = '  +0.54%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '89.42'
$ws.Range('E14').Value = '  +1.27%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6669'
$ws.Range('E15').Value = '  -0.53%  '

# Row 16
$ws.Range('D16').Value = '30.491.92'
$ws.Range('E16').Value = '  -0.74%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.58'
$ws.Range('E17').Value = '  +0.70%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007794'
$ws.Range('E18').Value = '  -1.41%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').Value = '  -0.07%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.440'
$ws.Range('E20').Value = '  +2.39%  '

# Row 21
$ws.Range('D21').Value = '2.141.86'
$ws.Range('E21').Value = '  -1.00%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.06%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '221.00'
$ws.Range('E23').Value = '  +12.74%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.214'
$ws.Range('E24').Value = '  -0.59%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.464'
$ws.Range('E25').Value = '  -1.70%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.00'
$ws.Range('E26').Value = '  +0.66%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.67'
$ws.Range('E27').Value = '  +0.42%  '

# Row 28
$ws.Range('E28').Value = '  +0.17%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.446'
$ws.Range('E29').Value = '  -1.47%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.331'
$ws.Range('E30').Value = '  -1.62%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09191'
$ws.Range('E31').Value = '  +0.74%  '

# Row 32
$ws.Range('E32').Value = '  +0.61%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05095'
$ws.Range('E33').Value = '  -2.80%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7531'
$ws.Range('E34').Value = '  +1.86%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.162'
$ws.Range('E35').Value = '  +4.42%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.702'
$ws.Range('E36').Value = '  -0.77%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01898'
$ws.Range('E37').Value = '  +4.34%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.654'
$ws.Range('E38').Value = '  -2.11%  '

# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.104'
$ws.Range('E39').Value = '  +1.37%  '

# Row 40
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9219'
$ws.Range('E40').Value = '  +0.10%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.069'
$ws.Range('E41').Value = '  +2.59%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '107.28'
$ws.Range('E42').Value = '  +0.58%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4363'
$ws.Range('E43').Value = '  -1.77%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.003'
$ws.Range('E44').Value = '  +0.29%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.668'
$ws.Range('E45').Value = '  +1.02%  '

# Row 46
$ws.Range('E46').Value = '  -2.64%  '

# Row 47
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.585'
$ws.Range('E47').Value = '  +10.88%  '

# Row 48
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.95'
$ws.Range('E48').Value = '  -11.42%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.946'
$ws.Range('E49').Value = '  -1.60%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.47'
$ws.Range('E50').Value = '  -2.43%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05702'
$ws.Range('E51').Value = '  -2.54%  '
